# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Change cell B11 on the active sheet from "R40" to the text value "1".
#
# Note: setting the cell's .Value directly to the string "1" would be
# auto-interpreted by Excel as a number, which would (a) store it as a
# numeric cell instead of a shared string and (b) churn the cell's style
# (quote-prefix gets applied). To keep the result as a genuine text value
# (shared string) while preserving the existing cell style, we write a
# formula that evaluates to the text "1" and then convert that formula to
# its static value via Copy / Paste-Special (values only), which is how
# Excel itself "flattens" a formula result into a literal without
# re-triggering general number inference on a plain string assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")

# Produce the literal text "1" as a formula result (forces text type).
$cell.Formula = "=TEXT(1,""0"")"

# Freeze the formula result into a static value, preserving the cell's
# existing style/formatting (xlPasteValues = -4163).
$cell.Copy()
$cell.PasteSpecial(-4163)

$excel.CutCopyMode = $false
